$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 10 (Friday 1.6.18): fill in the previously blank row.
$ws.Range("A10").Value = "Wk[12] Friday 1.6.18"
$ws.Range("B10").Value = "1400 - 1700 "
$ws.Range("C10").Value = 3

# Row 9 (Thursday 31.5.18): time text unchanged, description text is extended.
$ws.Range("D9").Value = "Testing simple structure (2 stages 1 queue)."

$ws.Range("D10").Value = "Implementing interactions between more complex structures."

# Keep the active selection where the user left off editing.
$ws.Range("D18").Select() | Out-Null
